# Refresh the cryptos list (prices / 1h volume, plus a few coins that
# changed rank and swapped rows) as published by the scheduled
# GitHub Actions job. All Price/Volume cells are stored as plain text
# in this sheet, so values that look numeric are written with a
# leading apostrophe to stop Excel from auto-converting them to
# Number (which would otherwise silently drop significant trailing
# zeros, e.g. "1.00" -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.760.26"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.347.19"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'323.02"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'102.46"
$ws.Range("E6").Value = "  -4.58%  "
$ws.Range("D7").Value = "'0.641"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.620"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").Value = "'39.93"
$ws.Range("E10").Value = "  -6.11%  "
$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "'8.48"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").Value = "'1.00"
$ws.Range("E13").Value = "  -3.60%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'16.07"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("D16").Value = "2.699.08"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "2.354.11"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "42.713.26"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "'7.81"
$ws.Range("E19").Value = "  +7.81%  "
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "'76.73"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Value = "'3.61"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "'266.24"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -7.07%  "
$ws.Range("D25").Value = "'9.79"
$ws.Range("E25").Value = "  +6.20%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'11.47"
$ws.Range("E27").Value = "  -4.51%  "
$ws.Range("D28").Value = "'22.95"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").Value = "'174.50"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "'6.25"
$ws.Range("E32").Value = "  +5.35%  "
$ws.Range("D33").Value = "'35.59"
$ws.Range("E33").Value = "  -9.26%  "
$ws.Range("D34").Value = "'0.0897"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "'0.133"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("D37").Value = "'4.56"
$ws.Range("E37").Value = "  -8.56%  "
$ws.Range("D38").Value = "'0.0359"
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").Value = "'3.79"
$ws.Range("E39").Value = "  -7.74%  "
$ws.Range("D40").Value = "'2.70"
$ws.Range("E40").Value = "  -4.25%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.236"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.48"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "'70.44"
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "'93.82"
$ws.Range("E44").Value = "  +22.96%  "
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'120.38"
$ws.Range("E46").Value = "  +6.79%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'5.56"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").Value = "'11.86"
$ws.Range("E48").Value = "  -6.40%  "
$ws.Range("D49").Value = "'9.14"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("E50").Value = "  -4.03%  "
$ws.Range("E51").Value = "  -0.28%  "
